# Generate Report for Archive
#
# 1) Shared-string text change: "Ready for handoff" -> "In Translation"
#    (appears in Overview!E2:F2,E3:F3 and in the Status column (C) of the
#    zh-cn / de-de sheets).
# 2) Narrow the "Status" / language-result columns (previously sized to fit
#    "Ready for handoff") now that the text is shorter:
#      - Overview sheet: columns E (zh-cn) and F (de-de)
#      - zh-cn sheet: column C (Status)
#      - de-de sheet: column C (Status)

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

# New column width to apply (character units). This lands the saved OOXML
# column width on the same pixel bucket as the target width.
$newColWidth = 12.58

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- 1. Replace the status text everywhere it occurs ---
# (Use Value2 for reads -- Value's getter in this host stringifies to the
#  member descriptor rather than the cell content; Value2 round-trips fine.)

if ($ws1.Range("E2").Value2 -eq $oldText) { $ws1.Range("E2").Value2 = $newText }
if ($ws1.Range("F2").Value2 -eq $oldText) { $ws1.Range("F2").Value2 = $newText }
if ($ws1.Range("E3").Value2 -eq $oldText) { $ws1.Range("E3").Value2 = $newText }
if ($ws1.Range("F3").Value2 -eq $oldText) { $ws1.Range("F3").Value2 = $newText }

if ($ws2.Range("C2").Value2 -eq $oldText) { $ws2.Range("C2").Value2 = $newText }
if ($ws2.Range("C3").Value2 -eq $oldText) { $ws2.Range("C3").Value2 = $newText }

if ($ws3.Range("C2").Value2 -eq $oldText) { $ws3.Range("C2").Value2 = $newText }
if ($ws3.Range("C3").Value2 -eq $oldText) { $ws3.Range("C3").Value2 = $newText }

# --- 2. Narrow the now-shorter columns ---

$ws1.Columns.Item(5).ColumnWidth = $newColWidth   # Overview column E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = $newColWidth   # Overview column F (de-de)
$ws2.Columns.Item(3).ColumnWidth = $newColWidth   # zh-cn column C (Status)
$ws3.Columns.Item(3).ColumnWidth = $newColWidth   # de-de column C (Status)
